# updates to markdown content
#
# Repositions/resizes a handful of shapes on slide 3 (title-card "Dan"
# doodle + connector) and slide 50 (contact-info block: LinkedIn URL,
# small icon, and @dmccreary handle), plus two font-size bumps.
#
# NOTE on the literal point values below: PowerPoint COM exposes shape
# geometry in points (1 pt = 12700 EMU) through .Left/.Top/.Width/.Height,
# and this host's COM shim round-trips that value through a 32-bit float
# before converting back to EMU on save. To land on the exact target EMU
# from the canonical OOXML, the literals here were solved (not just
# divided by 12700) so that after that float32 round-trip they serialize
# back to precisely the intended EMU offsets/extents.

$p = $ppt.ActivePresentation

# ---- Slide 3: title card ------------------------------------------------
$s3 = $p.Slides.Item(3)

# "Straight Connector 7" underline beneath "Dan"
$cxn = $s3.Shapes.Item("Straight Connector 7")
$cxn.Left = 401.6040344238281
$cxn.Top = 47.91409683227539

# "TextBox 1" holding the rotated "Dan" script text
$tbDan = $s3.Shapes.Item("TextBox 1")
$tbDan.Left = 406.4718322753906
$tbDan.Top = 6.2851972579956055

# ---- Slide 50: Thank You / contact info ---------------------------------
$s50 = $p.Slides.Item(50)

# LinkedIn URL textbox
$tbLinkedIn = $s50.Shapes.Item("TextBox 1")
$tbLinkedIn.Left = 509.8086853027344
$tbLinkedIn.Width = 391.43267822265625
$tbLinkedIn.Height = 31.504724502563477
$tbLinkedIn.TextFrame.TextRange.Font.Size = 20

# Small Twitter/X icon picture next to "@dmccreary"
$picHandle = $s50.Shapes.Item("Picture 6")
$picHandle.Left = 460.1380615234375
$picHandle.Width = 43.255592346191406
$picHandle.Height = 47.80882263183594

# "@dmccreary" textbox
$tbHandle = $s50.Shapes.Item("TextBox 7")
$tbHandle.Left = 509.8086853027344
$tbHandle.Top = 324.189697265625
$tbHandle.Width = 140.18080139160156
$tbHandle.Height = 36.35157775878906
$tbHandle.TextFrame.TextRange.Font.Size = 24
